# Insert a new header/description row at the very top of the "p1" sheet,
# pushing all existing rows down by one, and populate the new A1 cell
# with a "description" label (mirroring the header style used on the
# "index" sheet's A1, e.g. "header1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 1; existing content shifts to row 2+.
$ws.Rows.Item(1).Insert()

# New A1 acts as a row-label header, styled like the "header1" cells
# elsewhere in the workbook (e.g. sheet "index"!A1).
$excel.CutCopyMode = $false
$indexSheet = $wb.Worksheets.Item("index")
$indexSheet.Range("A1").Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 1).Value = "description"

# Restore the selection to B8 (matches the post-edit saved view state).
$ws.Range("B8").Select()
